$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2018-12-31 00:00:00"
$ws.Range("O2").Value = 194356845.08
$ws.Range("P2").Value = 9510.9783326012
$ws.Range("Q2").Value = 1991315077.4
$ws.Range("R2").Value = 97446.2954815812
$ws.Range("S2").Value = 446432829.47
$ws.Range("T2").Value = 21846.4802014219
$ws.Range("U2").Value = -149939510.81
$ws.Range("V2").Value = -7337.3872575867
$ws.Range("Y2").Value = 149964524.17
$ws.Range("Z2").Value = 7338.6113025895
$ws.Range("AA2").Value = -46069025.18
$ws.Range("AB2").Value = -2254.4176414815
$ws.Range("AC2").Value = 2043500.03
$ws.Range("AD2").Value = -97.8445420816
